$wb = $excel.ActiveWorkbook

# --- Clear the stale cached "Tax Base" growth values in column E on the
#     "Growth by Tax" sheet (they were left over from a previous revision of
#     the chart data and are no longer needed). ---
$wsGrowthByTax = $wb.Worksheets.Item("Growth by Tax")

$cellsToClear = @(
    "E9","E10","E11","E12","E13",
    "E16","E17","E18","E19","E20",
    "E23","E24","E25","E26","E27",
    "E30","E31","E32","E33","E34",
    "E37","E38","E39","E40","E41",
    "E44","E45","E46","E47","E48",
    "E51","E52","E53","E54","E55"
)

foreach ($addr in $cellsToClear) {
    $wsGrowthByTax.Range($addr).ClearContents()
}

# --- Update sheet view / selection state ---

# "Growth by Year" keeps its own selection updated to E9, without becoming
# the active sheet.
$wsGrowthByYear = $wb.Worksheets.Item("Growth by Year")
$wsGrowthByYear.Activate()
$wsGrowthByYear.Range("E9").Select()

# "Growth by Tax" becomes the active / selected tab with E9 selected.
$wsGrowthByTax.Activate()
$wsGrowthByTax.Range("E9").Select()

Write-Output "done"
